$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (summary) sheet
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(3)
$q1 = $wb.Worksheets.Add($total)
$q1.Name = "2022-Q1"

# Reference sheet that already carries the header/"index" style (s="2")
$styleRef = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Helper-ish inline pattern: write text-like numeric strings without letting
# Excel auto-convert them to real numbers (and without leaving a stray
# custom number-format behind once we are done).
# ---------------------------------------------------------------------------

# -- header row (B1:H1), styled like the other sheets' headers --------------
$styleRef.Range("B1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# -- data row 2 ---------------------------------------------------------------
$styleRef.Range("A2").Copy()
$q1.Range("A2").PasteSpecial(-4122)
$q1.Range("A2").Value = 0

$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "005269"
$q1.Range("B2").Style = "Normal"

$q1.Range("C2").Value = "华泰柏瑞港股通量化灵活配置混合"

$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "0.33"
$q1.Range("D2").Style = "Normal"

$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "37.77"
$q1.Range("E2").Style = "Normal"

$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "1.06"
$q1.Range("F2").Style = "Normal"

$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "0.0035"
$q1.Range("G2").Style = "Normal"

$q1.Range("H2").Value = 4

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new top data row for 2022-Q1 and push
#    the existing 2021-Q4 / 2021-Q3 rows down by one.
# ---------------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item(4)

# push 2021-Q3 (currently row 3) down to row 4
$totalWs.Range("A4").Value = 2
$totalWs.Range("B4").Value = "2021-Q3"
$totalWs.Range("C4").Value = 3
$totalWs.Range("D4").Value = 0.18

# push 2021-Q4 (currently row 2) down to row 3
$totalWs.Range("A3").Value = 1
$totalWs.Range("B3").Value = "2021-Q4"
$totalWs.Range("C3").Value = 2
$totalWs.Range("D3").Value = 0.22

# new top row: 2022-Q1
$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("C2").Value = 1
$totalWs.Range("D2").Value = 0

# keep the index-column style (s="2") consistent on every row
$totalWs.Range("A2").Copy()
$totalWs.Range("A3:A4").PasteSpecial(-4122)
$totalWs.Range("A2").Value = 0
$totalWs.Range("A3").Value = 1
$totalWs.Range("A4").Value = 2
